$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 77.75
$ws.Range("C2").Value = 'The applicant has worked on a variety of projects that demonstrate strong technical skills and experience in the relevant areas. The project ''LLMGuard'' showcases the use of NLP models like GPT-2 and BERT, as well as the development of a user-friendly UI and backend, highlighting expertise in NLP, Python, and ReactJS. Additionally, the project ''CloudPhysician''s Vital Extraction Challenge'' exhibits proficiency in Computer Vision and Image Processing, essential for the role. Although the applicant lacks direct experience with PyTorch, their experience with Tensorflow and Python demonstrates a strong foundation in deep learning frameworks. The applicant''s wide range of technical skills and hands-on experience make them well-suited for the NLP Engineer role.'
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 'The applicant''s answers reflect a strong willingness to work in Japan, emphasizing a desire to learn from Japanese work culture and citing personal alignment with Japanese values like hard work and dedication. The applicant also expresses a clear understanding of their strengths and areas for improvement, showcasing strong self-awareness and a commitment to continuous self-improvement. Their responses demonstrate adaptability, strong communication skills, and a positive attitude, indicating a good fit for the company''s cultural values and the willingness to work in Japan.'

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 66.82
$ws.Range("C3").Value = 'The applicant has a good score based on the projects he has worked on, which showcase his proficiency in a variety of technologies such as ReactJS, Flutter, WebRTC, and NodeJS. Although the job requires skills in NLP, Pytorch, and Computer Vision, the applicant''s strong technical background and experience in building web and mobile applications demonstrate his ability to learn and adapt to new technologies. Additionally, the applicant''s experience in developing real-time communication applications using WebSockets and WebRTC API aligns with the problem-solving and adaptability skills required for the role.'
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 'The applicant has displayed a strong willingness to work in Japan, citing an interest in the work culture and expressing readiness to adapt to a new environment. The applicant''s desire to work on cutting-edge AI/ML technologies, coupled with the openness to learn a new language, demonstrates a strong motivation and adaptability. Furthermore, the applicant''s emphasis on teamwork and collaborative problem-solving align with the soft skills sought by the company, indicating a good fit for the company''s culture.'

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 89.94
$ws.Range("C4").Value = 'The applicant has extensive experience in projects related to NLP, Pytorch, and Computer Vision, which are directly in line with the company''s requirements. The applicant has demonstrated hands-on experience in developing and improving models for NLP and utilizing Pytorch and Computer Vision in various projects. Although some specific skills like NLP are not explicitly mentioned in the projects, the applicant''s overall technical expertise and the demonstrated ability to work with similar technologies make them well-suited for the role. The diverse range of projects showcases the applicant''s adaptability and problem-solving skills, which are crucial for the position. However, the lack of explicit NLP project experience could be a minor concern, but the applicant''s strong foundation in related domains compensates for it.'
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 'The applicant has expressed a keen interest in working in Japan and has demonstrated adaptability to new cultures. They have highlighted the opportunity to learn a new language and expressed a strong desire to be a part of Japan''s work culture. Their willingness to work in Japan and openness to adaptability align with the company''s requirements for a candidate who is receptive to new environments. Additionally, the applicant has emphasized their ability to work well in teams and their commitment to contribute to cutting-edge technologies, reflecting qualities such as teamwork, problem-solving, and adaptability, which are desirable traits for the company. The applicant''s proactive approach to addressing weaknesses and leveraging their unique background further exhibits a willingness to grow and succeed in a new professional setting.'

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 60.16
$ws.Range("C5").Value = 'The applicant has demonstrated strong skills in JavaScript, HTML, and CSS through projects like ''Video Conferencing Project'' and ''TEDxIITGuwahati Website''. Although these skills are not directly related to the required skills of NLP, Pytorch, Computer Vision, and Python, they showcase the applicant''s proficiency in web development, which could be valuable for frontend development or interface design within the NLP field. The applicant''s API integration skills also demonstrate an understanding of integrating different systems, which could be transferable to the NLP domain. While the applicant''s projects do not directly align with the required skills, they exhibit a strong foundation in programming and web development, which can be beneficial for certain aspects of the role.'
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 'Based on the applicant''s answers, it is evident that the applicant is deeply interested in AI and is willing to work in Japan. The applicant has expressed a strong fascination for the work culture in Japan and has shown a willingness to adapt to a new culture and learn a new language. Additionally, the applicant has demonstrated good communication skills and a positive attitude towards working in a team. Overall, the applicant''s responses indicate a strong willingness to work in Japan and align with the company''s soft skills expectations, hence receiving a score of 4 out of 5.'
